# Commit by somnath as on 16/06/2021
#
# Changes applied:
#  1. Add a new worksheet "OrderWithRewards" at the end of the workbook with
#     two rows of data (Reward Item / Chips and Salsa).
#  2. Swap the order of the two locations on the "Locations" sheet
#     (Delray Beach now first, Addison second) and move the saved
#     selection there to D15.
#  3. Swap the paired Store Name values on the "LoggedInOrder" sheet
#     (coit-road / delray-beach) to match the new Locations order, and
#     make "LoggedInOrder" the active / selected sheet (selection D4)
#     instead of "UpdateMyAccount".
$wb = $excel.ActiveWorkbook

# --- 1. New "OrderWithRewards" sheet, inserted after the last tab ---------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$rewardsSheet = $wb.Worksheets.Add($null, $lastSheet)
$rewardsSheet.Name = "OrderWithRewards"
$rewardsSheet.Range("A1").Value = "Reward Item"
$rewardsSheet.Range("A2").Value = "Chips and Salsa"
$rewardsSheet.Columns.Item(1).AutoFit() | Out-Null

# --- 2. Swap Locations order (sheet "Locations") ---------------------------
$locSheet = $wb.Worksheets.Item("Locations")
$oldA2 = $locSheet.Range("A2").Value()
$oldA3 = $locSheet.Range("A3").Value()
$locSheet.Range("A2").Value = $oldA3
$locSheet.Range("A3").Value = $oldA2
$locSheet.Range("D15").Select() | Out-Null

# --- 3. Swap Store Name values on "LoggedInOrder" and make it active ------
# A2 held "delray-beach" formatted as General; A3 held "coit-road" formatted
# as Text ("@", style index 3). The values AND their formatting both swap
# together (as if the two cells' contents were exchanged), so after the
# edit A2 carries the Text format and A3 reverts to the default format.
$orderSheet = $wb.Worksheets.Item("LoggedInOrder")
$a2Cell = $orderSheet.Range("A2")
$a3Cell = $orderSheet.Range("A3")
$oldA2b = $a2Cell.Value()
$oldA3b = $a3Cell.Value()

$a2Cell.Value = $oldA3b
$a3Cell.Value = $oldA2b

$a2Cell.NumberFormat = "@"
$a3Cell.ClearFormats() | Out-Null

$orderSheet.Activate()
$orderSheet.Range("D4").Select() | Out-Null
